$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each new row: copy formatting from the last existing data row (323),
# then overwrite per-cell values. Columns B and C stay blank (as in the
# source row); A and L repeat constant values also found on row 323.

# ---- Row 324 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A324:M324").PasteSpecial(-4122)
$ws.Range("A324").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D324").PasteSpecial(-4163)
$ws.Range("E324").Value = "Event month"
$ws.Range("F324").Value = "1: 1414"
$ws.Range("G324").Value = "1: 1421"
$ws.Range("H324").Value = 0
$ws.Range("I324").Value = "December"
$ws.Range("J324").Value = 8
$ws.Range("K324").Value = 0.016161
$ws.Range("L324").Value = "Sonia"
$ws.Range("M324").Value = "11/14/18 10:06:00"

# ---- Row 325 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A325:M325").PasteSpecial(-4122)
$ws.Range("A325").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D325").PasteSpecial(-4163)
$ws.Range("E325").Value = "Event month"
$ws.Range("F325").Value = "1: 1431"
$ws.Range("G325").Value = "1: 1433"
$ws.Range("H325").Value = 0
$ws.Range("I325").Value = "May"
$ws.Range("J325").Value = 3
$ws.Range("K325").Value = 0.00606
$ws.Range("L325").Value = "Sonia"
$ws.Range("M325").Value = "11/14/18 10:06:00"

# ---- Row 326 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A326:M326").PasteSpecial(-4122)
$ws.Range("A326").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D326").PasteSpecial(-4163)
$ws.Range("E326").Value = "Event year"
$ws.Range("F326").Value = "1: 1423"
$ws.Range("G326").Value = "1: 1426"
$ws.Range("H326").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2004,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I326").PasteSpecial(-4163)
$ws.Range("J326").Value = 4
$ws.Range("K326").Value = 0.008081
$ws.Range("L326").Value = "Sonia"
$ws.Range("M326").Value = "11/14/18 10:06:00"

# ---- Row 327 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A327:M327").PasteSpecial(-4122)
$ws.Range("A327").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D327").PasteSpecial(-4163)
$ws.Range("E327").Value = "Event year"
$ws.Range("F327").Value = "1: 1435"
$ws.Range("G327").Value = "1: 1438"
$ws.Range("H327").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2005,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I327").PasteSpecial(-4163)
$ws.Range("J327").Value = 4
$ws.Range("K327").Value = 0.008081
$ws.Range("L327").Value = "Sonia"
$ws.Range("M327").Value = "11/14/18 10:06:00"

# ---- Row 328 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A328:M328").PasteSpecial(-4122)
$ws.Range("A328").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D328").PasteSpecial(-4163)
$ws.Range("E328").Value = "A"
$ws.Range("F328").Value = "1: 1414"
$ws.Range("G328").Value = "1: 1421"
$ws.Range("H328").Value = 0
$ws.Range("I328").Value = "December"
$ws.Range("J328").Value = 8
$ws.Range("K328").Value = 0.016161
$ws.Range("L328").Value = "Sonia"
$ws.Range("M328").Value = "11/14/18 10:06:00"

# ---- Row 329 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A329:M329").PasteSpecial(-4122)
$ws.Range("A329").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D329").PasteSpecial(-4163)
$ws.Range("E329").Value = "A"
$ws.Range("F329").Value = "1: 1423"
$ws.Range("G329").Value = "1: 1426"
$ws.Range("H329").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2004,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I329").PasteSpecial(-4163)
$ws.Range("J329").Value = 4
$ws.Range("K329").Value = 0.008081
$ws.Range("L329").Value = "Sonia"
$ws.Range("M329").Value = "11/14/18 10:06:00"

# ---- Row 330 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A330:M330").PasteSpecial(-4122)
$ws.Range("A330").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D330").PasteSpecial(-4163)
$ws.Range("E330").Value = "B"
$ws.Range("F330").Value = "1: 1431"
$ws.Range("G330").Value = "1: 1433"
$ws.Range("H330").Value = 0
$ws.Range("I330").Value = "May"
$ws.Range("J330").Value = 3
$ws.Range("K330").Value = 0.00606
$ws.Range("L330").Value = "Sonia"
$ws.Range("M330").Value = "11/14/18 10:07:00"

# ---- Row 331 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A331:M331").PasteSpecial(-4122)
$ws.Range("A331").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(10703,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D331").PasteSpecial(-4163)
$ws.Range("E331").Value = "B"
$ws.Range("F331").Value = "1: 1435"
$ws.Range("G331").Value = "1: 1438"
$ws.Range("H331").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2005,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I331").PasteSpecial(-4163)
$ws.Range("J331").Value = 4
$ws.Range("K331").Value = 0.008081
$ws.Range("L331").Value = "Sonia"
$ws.Range("M331").Value = "11/14/18 10:07:00"

# ---- Row 332 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A332:M332").PasteSpecial(-4122)
$ws.Range("A332").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(11202,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D332").PasteSpecial(-4163)
$ws.Range("E332").Value = "Event month"
$ws.Range("F332").Value = "1: 506"
$ws.Range("G332").Value = "1: 513"
$ws.Range("H332").Value = 0
$ws.Range("I332").Value = "December"
$ws.Range("J332").Value = 8
$ws.Range("K332").Value = 0.041216
$ws.Range("L332").Value = "Sonia"
$ws.Range("M332").Value = "11/14/18 10:07:00"

# ---- Row 333 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A333:M333").PasteSpecial(-4122)
$ws.Range("A333").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(11202,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D333").PasteSpecial(-4163)
$ws.Range("E333").Value = "Event month"
$ws.Range("F333").Value = "1: 523"
$ws.Range("G333").Value = "1: 530"
$ws.Range("H333").Value = 0
$ws.Range("I333").Value = "December"
$ws.Range("J333").Value = 8
$ws.Range("K333").Value = 0.041216
$ws.Range("L333").Value = "Sonia"
$ws.Range("M333").Value = "11/14/18 10:07:00"

# ---- Row 334 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A334:M334").PasteSpecial(-4122)
$ws.Range("A334").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(11202,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D334").PasteSpecial(-4163)
$ws.Range("E334").Value = "Event year"
$ws.Range("F334").Value = "1: 515"
$ws.Range("G334").Value = "1: 518"
$ws.Range("H334").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2011,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I334").PasteSpecial(-4163)
$ws.Range("J334").Value = 4
$ws.Range("K334").Value = 0.020608
$ws.Range("L334").Value = "Sonia"
$ws.Range("M334").Value = "11/14/18 10:07:00"

# ---- Row 335 ----
$ws.Range("A323:M323").Copy()
$ws.Range("A335:M335").PasteSpecial(-4122)
$ws.Range("A335").Value = "●"
$ws.Range("ZZ1").Formula = "=TEXT(11202,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D335").PasteSpecial(-4163)
$ws.Range("E335").Value = "Event year"
$ws.Range("F335").Value = "1: 532"
$ws.Range("G335").Value = "1: 535"
$ws.Range("H335").Value = 0
$ws.Range("ZZ1").Formula = "=TEXT(2012,""0"")"
$ws.Range("ZZ1").Copy()
$ws.Range("I335").PasteSpecial(-4163)
$ws.Range("J335").Value = 4
$ws.Range("K335").Value = 0.020608
$ws.Range("L335").Value = "Sonia"
$ws.Range("M335").Value = "11/14/18 10:08:00"

# cleanup scratch cell
$ws.Range("ZZ1").Clear()
$ws.Application.CutCopyMode = $false
